# Update rules for prospec
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that are being removed (the ">2500.00" and ">6" rule
# rows), which shifts the remaining rows (old 6,7,8) up into rows 4,5,6.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Row 4 (was old row 6): FedEx NextDay Styrofoam Box rule - unchanged values,
# just confirm content lines up with spec after the shift.
$ws.Range("A4").Value = "`$(ProSpec) <=  `$in"
$ws.Range("B4").Value = "`$in >= `$(ice-packs)"
$ws.Range("C4").Value = "<7"
$ws.Range("D4").Value = "<2500.01"
$ws.Range("E4").Value = "FedEx NextDay Styrofoam Box"
$ws.Range("F4").Value = "ice packs/blue ice"
$ws.Range("G4").Value = 2000
$ws.Range("H4").Value = 8500
$ws.Range("I4").Value = "N/A"
$ws.Range("J4").Value = "N/A"
$ws.Range("K4").Value = "N/A"
$ws.Range("L4").Value = "N/A"

# Row 5 (was old row 7): FedEx International Express Small Box rule - now
# also specifies the total-value tier and a handling note for lyophilized
# freeze-dry product.
$ws.Range("A5").Value = "`$(ProSpec) <=  `$in"
$ws.Range("B5").Value = "`$(ice-packs) intersection `$in #=0"
$ws.Range("C5").Value = "<7"
$ws.Range("D5").Value = "<2500.01"
$ws.Range("E5").Value = "FedEx International Express Small Box"
$ws.Range("F5").Value = "lyophilized freeze-dry"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 4500
$ws.Range("I5").Value = 5500
$ws.Range("J5").Value = 6500
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 8000

# Row 6 (was old row 8): FedEx NextDay Small Box rule - same new
# total-value tier and lyophilized freeze-dry handling note.
$ws.Range("A6").Value = "`$(ProSpec) <=  `$in"
$ws.Range("B6").Value = "`$(ice-packs) intersection `$in #=0"
$ws.Range("C6").Value = "<7"
$ws.Range("D6").Value = "<2500.01"
$ws.Range("E6").Value = "FedEx NextDay Small Box"
$ws.Range("F6").Value = "lyophilized freeze-dry"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 5500
$ws.Range("I6").Value = "N/A"
$ws.Range("J6").Value = "N/A"
$ws.Range("K6").Value = "N/A"
$ws.Range("L6").Value = "N/A"

# Update the view: the sheet now ends at row 6, scrolled right to column F,
# and the full used range A1:XFD6 is selected (rather than a stray G12).
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("A1:XFD6").Select()
